$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Performance ..." bullet so we
# can add the new bullet right after it (and before the trailing empty
# "List Paragraph" paragraph).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Performance*dispatchers?*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Performance' bullet paragraph"
}

$r = $target.Range
$r.InsertParagraphAfter()

# The freshly inserted paragraph now sits right after $target; it
# inherited $target's paragraph style + list numbering, matching the
# surrounding bullets.
$newPara = $target.Next()
$newRange = $newPara.Range
$newRange.InsertBefore("Fault-tolerance / reliability – have a visual display of all messages being sent on the dead letters actor channel – check them regularly.")

Write-Output "inserted new task paragraph"
